# Applies the "Added Figure S2 (spec curve)" refresh to TableS14.xlsx:
# the underlying TableS14.csv was re-generated (new spec-curve numbers),
# Power Query re-imported it as a new query/table generation
# "TableS14 (2)" -> "TableS14 (3)", and the displayed values were updated
# to match the refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TableS14 (2)")

# --- Rename the worksheet (cascades to the sheet's defined name formula) ---
$ws.Name = "TableS14 (3)"

# --- Rename the query/auto-filter table backing the sheet ---
$lo = $ws.ListObjects.Item(1)
$lo.Name = "TableS14__3"

# --- Rename the ExternalData defined name (2 -> 3) ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*ExternalData_2*") {
        $n.Name = "ExternalData_3"
    }
}

# --- Update the refreshed data values (rows 4-6, columns B-D) ---
# Force Text storage (matching the Power Query "type text" transform used
# for every column) so values like " 0.026" keep their leading space and
# are not coerced to numbers.
$dataRng = $ws.Range("B4:D6")
$dataRng.NumberFormat = "@"

$ws.Cells.Item(4, 2).Value = "-0.057"
$ws.Cells.Item(4, 3).Value = " 0.026"
$ws.Cells.Item(4, 4).Value = "-0.083"

$ws.Cells.Item(5, 2).Value = "-0.045"
$ws.Cells.Item(5, 3).Value = " 0.012"
$ws.Cells.Item(5, 4).Value = "-0.057"

$ws.Cells.Item(6, 4).Value = "68.7"

# --- Column widths tightened slightly after the refresh (best effort) ---
$ws.Columns.Item(1).ColumnWidth = 31.3
$ws.Columns.Item(2).ColumnWidth = 18.9
$ws.Columns.Item(3).ColumnWidth = 18.5
$ws.Columns.Item(4).ColumnWidth = 10.7

# --- Update selection to match the refreshed table ---
[void]$ws.Range("B4:D6").Select()
